$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Find-ParagraphIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd("`r", "`n", [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# Inserts a brand-new "ListParagraph"-styled paragraph (Times New Roman 12pt,
# spacing-after 0, bulleted/numbered with the given numId) right after the
# paragraph at $afterIndex. $runsXml is the raw <w:r>...</w:r> markup (one or
# more runs) that becomes the paragraph's content.
function Insert-ListParagraphAfter($doc, $afterIndex, $numId, $runsXml) {
    $p = $doc.Paragraphs.Item($afterIndex)
    $r = $p.Range
    $r.InsertParagraphAfter()
    $newp = $doc.Paragraphs.Item($afterIndex + 1)
    $newr = $newp.Range

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p>' +
           '<w:pPr>' +
           '<w:pStyle w:val="ListParagraph"/>' +
           '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr>' +
           '<w:spacing w:after="0"/>' +
           '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
           '</w:pPr>' +
           $runsXml +
           '</w:p>' +
           '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $newr.InsertXML($xml)
    return $newp
}

function RunXml($text, [bool]$preserve) {
    if ($preserve) {
        $space = ' xml:space="preserve"'
    } else {
        $space = ''
    }
    return '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t' + $space + '>' + $text + '</w:t></w:r>'
}

# ---------------------------------------------------------------------------
# 1) "Upload profile pic. ?" - new bullet (numId 2) after "Suspend account. ?"
# ---------------------------------------------------------------------------

$idx1 = Find-ParagraphIndexByText $d "Suspend account. ?"
if ($idx1 -lt 0) { throw "Could not locate 'Suspend account. ?' paragraph" }
$runs1 = RunXml "Upload profile pic. ?" $false
Insert-ListParagraphAfter $d $idx1 "2" $runs1 | Out-Null

# ---------------------------------------------------------------------------
# 2) "Upload pictures. ?" (3 runs) - new bullet (numId 4) after
#    "Notify Producer/Publisher. ?"
# ---------------------------------------------------------------------------

$idx2 = Find-ParagraphIndexByText $d "Notify Producer/Publisher. ?"
if ($idx2 -lt 0) { throw "Could not locate 'Notify Producer/Publisher. ?' paragraph" }
$runs2 = (RunXml "Upload pictures" $false) + (RunXml "." $false) + (RunXml " ?" $true)
Insert-ListParagraphAfter $d $idx2 "4" $runs2 | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the stray <w:lastRenderedPageBreak/> in front of "Design:"
# ---------------------------------------------------------------------------

$null = $d.Content.Find.Execute("Design:", $true, $false, $false, $false, $false, $true, 1, $false, "Design:", 2)

# ---------------------------------------------------------------------------
# 4) "Upload pictures. " - new bullet (numId 8) after "Style interior."
# ---------------------------------------------------------------------------

$idx4 = Find-ParagraphIndexByText $d "Style interior."
if ($idx4 -lt 0) { throw "Could not locate 'Style interior.' paragraph" }
$runs4 = RunXml "Upload pictures. " $true
Insert-ListParagraphAfter $d $idx4 "8" $runs4 | Out-Null
